$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.1648366424899244
    "C2" = 0.6565107402614302
    "D2" = 0.6208915960224134
    "E2" = 0.7879667480436046
    "F2" = 0.7996194615330313

    "B3" = 0.06925299554247144
    "C3" = 0.6317528524474937
    "D3" = 0.5154239904830427
    "E3" = 0.7179303521115699
    "F3" = 0.7437609478291553

    "B4" = 0.1042866910554764
    "C4" = 0.5871660518276492
    "D4" = 0.4576314133407699
    "E4" = 0.6764845994852875
    "F4" = 0.6981187964102095

    "B5" = 0.1411991018086509
    "C5" = 0.5954485264947712
    "D5" = 0.500182751991249
    "E5" = 0.7072359945529136
    "F5" = 0.7268219329406943

    "B6" = 0.3268981278322622
    "C6" = 0.6566970789475719
    "D6" = 0.6666380539392567
    "E6" = 0.8164790590941428
    "F6" = 0.7886526259380041

    "B7" = 0.2872108118487342
    "C7" = 0.7140200429738547
    "D7" = 0.7624163934427893
    "E7" = 0.8731645855408872
    "F7" = 0.8745954126766147

    "B8" = 0.2887732378318851
    "C8" = 0.6550352748098781
    "D8" = 0.6159387406405581
    "E8" = 0.7848176480180336
    "F8" = 0.7994113517477577

    "B9" = -0.05480200069631858
    "C9" = 0.1889080360497501
    "D9" = 0.05129867587295883
    "E9" = 0.2264921099574085
    "F9" = 0.2691526052055958

    "B10" = -0.1770340777201527
    "C10" = 0.1770340777201527
    "D10" = 0.03134106467422507
    "E10" = 0.1770340777201527
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
